$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Column D holds price values formatted as plain text (e.g. "22.939.48" which
# is not a valid number). Force text format so Excel's COM layer does not
# silently re-interpret numeric-looking strings (e.g. "299.16") as numbers.
$ws.Range("D2:D51").NumberFormat = "@"

$ws.Range("D2").Value = "22.945.88"
$ws.Range("E2").Value = "  -0.97%  "
$ws.Range("D3").Value = "1.573.32"
$ws.Range("E3").Value = "  -2.07%  "
$ws.Range("E4").Value = "  +0.15%  "
$ws.Range("E5").Value = "  +0.11%  "
$ws.Range("D6").Value = "299.16"
$ws.Range("E6").Value = "  -1.03%  "
$ws.Range("D7").Value = "0.3748"
$ws.Range("E7").Value = "  -0.17%  "
$ws.Range("D8").Value = "0.3547"
$ws.Range("E8").Value = "  -2.29%  "
$ws.Range("D9").Value = "49.90"
$ws.Range("E9").Value = "  +2.67%  "
$ws.Range("D10").Value = "1.002"
$ws.Range("E10").Value = "  -0.01%  "
$ws.Range("D11").Value = "1.209"
$ws.Range("E11").Value = "  -4.15%  "
$ws.Range("D12").Value = "0.07948"
$ws.Range("E12").Value = "  -1.22%  "
$ws.Range("D13").Value = "21.72"
$ws.Range("E13").Value = "  -4.98%  "
$ws.Range("D14").Value = "6.381"
$ws.Range("E14").Value = "  -2.40%  "
$ws.Range("D15").Value = "7.266"
$ws.Range("E15").Value = "  -5.08%  "
$ws.Range("D16").Value = "0.00001219"
$ws.Range("E16").Value = "  -3.43%  "
$ws.Range("D17").Value = "1.576.33"
$ws.Range("E17").Value = "  -1.80%  "
$ws.Range("D18").Value = "91.78"
$ws.Range("E18").Value = "  +0.53%  "
$ws.Range("D19").Value = "0.06728"
$ws.Range("E19").Value = "  -0.79%  "
$ws.Range("D20").Value = "17.60"
$ws.Range("E20").Value = "  -3.78%  "
$ws.Range("E21").Value = "  +0.19%  "
$ws.Range("D22").Value = "6.341"
$ws.Range("E22").Value = "  -3.31%  "
$ws.Range("D23").Value = "22.973.39"
$ws.Range("E23").Value = "  -0.94%  "
$ws.Range("E24").Value = "  -3.71%  "
$ws.Range("E25").Value = "  +1.08%  "
$ws.Range("D26").Value = "2.787"
$ws.Range("E26").Value = "  -3.79%  "
$ws.Range("D27").Value = "20.53"
$ws.Range("E27").Value = "  -2.31%  "
$ws.Range("D28").Value = "147.32"
$ws.Range("E28").Value = "  -1.79%  "
$ws.Range("D29").Value = "5.160"
$ws.Range("E29").Value = "  -2.04%  "
$ws.Range("D30").Value = "130.76"
$ws.Range("E30").Value = "  -1.09%  "
$ws.Range("D31").Value = "2.329"
$ws.Range("E31").Value = "  -3.54%  "
$ws.Range("D32").Value = "6.510"
$ws.Range("E32").Value = "  -3.08%  "
$ws.Range("D33").Value = "1.749.80"
$ws.Range("E33").Value = "  -1.99%  "
$ws.Range("D34").Value = "0.9283"
$ws.Range("E34").Value = "  -4.13%  "
$ws.Range("D35").Value = "0.07330"
$ws.Range("E35").Value = "  -4.93%  "
$ws.Range("D36").Value = "0.08739"
$ws.Range("E36").Value = "  -1.15%  "
$ws.Range("D37").Value = "9.904"
$ws.Range("E37").Value = "  -1.50%  "
$ws.Range("D38").Value = "0.02624"
$ws.Range("E38").Value = "  -5.12%  "
$ws.Range("D39").Value = "0.2455"
$ws.Range("E39").Value = "  -3.15%  "
$ws.Range("D40").Value = "5.947"
$ws.Range("E40").Value = "  -4.28%  "
$ws.Range("D41").Value = "1.341"
$ws.Range("E41").Value = "  -3.74%  "
$ws.Range("D42").Value = "0.6851"
$ws.Range("E42").Value = "  -3.74%  "
$ws.Range("D43").Value = "11.77"
$ws.Range("E43").Value = "  -7.41%  "
$ws.Range("D44").Value = "14.67"
$ws.Range("E44").Value = "  -8.06%  "
$ws.Range("D45").Value = "1.001"
$ws.Range("E45").Value = "  +0.06%  "
$ws.Range("D46").Value = "0.6301"
$ws.Range("E46").Value = "  -3.97%  "
$ws.Range("D47").Value = "3.968"
$ws.Range("E47").Value = "  -0.23%  "
$ws.Range("D48").Value = "2.238"
$ws.Range("E48").Value = "  -2.02%  "
$ws.Range("B49").Value = "Cronos"
$ws.Range("C49").Value = "https://coinranking.com/coin/65PHZTpmE55b+cronos-cro"
$ws.Range("D49").Value = "0.07834"
$ws.Range("E49").Value = "  -1.91%  "
$ws.Range("B50").Value = "Quant"
$ws.Range("C50").Value = "https://coinranking.com/coin/bauj_21eYVwso+quant-qnt"
$ws.Range("D50").Value = "129.36"
$ws.Range("E50").Value = "  -1.42%  "
$ws.Range("D51").Value = "1.178"
$ws.Range("E51").Value = "  +1.27%  "
